$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    } elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "system, System, backup@backdoor.com"
    } elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value2 = "System, backup@backdoor.com"
    }
}
